$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "No se ha introducido ninguna URL"
$ws.Range("A18").Value = "No hay conexión a internet"
$ws.Range("A19").Value = "La URL no es de YouTube"

$ws.Range("B17").Value = "No URL entered"
$ws.Range("B18").Value = "No internet connection"
$ws.Range("B19").Value = "The URL is not from YouTube"

$ws.Range("A2:B2").Copy()
$ws.Range("A17:B19").PasteSpecial(-4122)

$ws.Rows("17:19").RowHeight = 15.75

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B19").Select()
